$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.341.41'
$ws.Cells.Item(2, 5).Value = '  +0.13%  '

$ws.Cells.Item(3, 4).Value = '3.408.39'
$ws.Cells.Item(3, 5).Value = '  -0.36%  '

$ws.Cells.Item(4, 5).Value = '  -0.01%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '569.41'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.43%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '156.56'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.80%  '

$ws.Cells.Item(7, 5).Value = '  +8.85%  '

$ws.Cells.Item(8, 5).Value = '  -0.15%  '

$ws.Cells.Item(9, 4).Value = '3.417.28'
$ws.Cells.Item(9, 5).Value = '  -0.07%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.13'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.39%  '

$ws.Cells.Item(11, 5).Value = '  -1.42%  '

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.440'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.09%  '

$ws.Cells.Item(13, 4).Value = '3.994.79'
$ws.Cells.Item(13, 5).Value = '  -0.40%  '

$ws.Cells.Item(14, 5).Value = '  -0.02%  '

$ws.Cells.Item(15, 5).Value = '  -2.17%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.55'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -1.40%  '

$ws.Cells.Item(17, 4).Value = '64.381.68'
$ws.Cells.Item(17, 5).Value = '  +0.00%  '

$ws.Cells.Item(18, 4).Value = '3.404.01'
$ws.Cells.Item(18, 5).Value = '  -0.96%  '

$ws.Cells.Item(19, 5).Value = '  -0.09%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.83'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.10%  '

$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '377.75'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -1.27%  '

$ws.Cells.Item(22, 5).Value = '  -1.21%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.544'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.77%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.70'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.62%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000118'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -2.47%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.35'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +7.82%  '

$ws.Cells.Item(28, 5).Value = '  -0.18%  '

$ws.Cells.Item(29, 5).Value = '  -0.07%  '

$ws.Cells.Item(30, 5).Value = '  +4.22%  '

$ws.Cells.Item(31, 5).Value = '  -0.15%  '

$ws.Cells.Item(32, 5).Value = '  -1.46%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.08'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -1.59%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.14'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.29%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.60'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +7.99%  '

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '159.49'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -1.81%  '

$ws.Cells.Item(37, 5).Value = '  +1.28%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.00'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +7.36%  '

$ws.Cells.Item(39, 5).Value = '  +0.08%  '

$ws.Cells.Item(40, 4).Value = '2.878.39'
$ws.Cells.Item(40, 5).Value = '  -4.64%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.61'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.53%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.28'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -2.67%  '

$ws.Cells.Item(43, 5).Value = '  +0.25%  '

$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0316'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -0.17%  '

$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '25.98'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +5.72%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.770'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.52%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '321.09'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +7.71%  '

$ws.Cells.Item(48, 5).Value = '  +0.58%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.110'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +3.40%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.18'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.84%  '

$ws.Cells.Item(51, 2).Value = 'Cosmos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.56'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.35%  '
